$wb = $excel.ActiveWorkbook

# --- Sheet ALC: market price refresh ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 66943.46000000001
$ws.Range("I98").Value = 66943.46000000001
$ws.Range("K98").Value = 66943.46000000001
$ws.Range("M98").Value = -65445.46000000001
$ws.Range("H112").Value = 4306.8184
$ws.Range("J112").Value = 4408
$ws.Range("L112").Value = 13224
$ws.Range("N112").Value = -15440
$ws.Range("H122").Value = 66943.46000000001
$ws.Range("I122").Value = 66943.46000000001
$ws.Range("K122").Value = 200830.38
$ws.Range("M122").Value = -198380.38
$ws.Range("H135").Value = 13890849
$ws.Range("J135").Value = 27779932
$ws.Range("L135").Value = 250019388
$ws.Range("N135").Value = -250024458
$ws.Range("H138").Value = 3394.0425
$ws.Range("J138").Value = 3831.0881
$ws.Range("L138").Value = 11493.2643
$ws.Range("N138").Value = -21773.2643

# --- Sheet ARM: market price refresh ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2367.4614
$ws.Range("I122").Value = 1788.875
$ws.Range("J122").Value = 3293.2
$ws.Range("K122").Value = 5366.625
$ws.Range("L122").Value = 9879.599999999999
$ws.Range("M122").Value = -2916.625
$ws.Range("N122").Value = -14779.6
$ws.Range("H132").Value = 23816666
$ws.Range("I132").Value = 7550.9473
$ws.Range("K132").Value = 22652.8419
$ws.Range("M132").Value = -20122.8419

# --- Sheet BSM: market price refresh ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 13743.286
$ws.Range("I75").Value = 6867.1665
$ws.Range("K75").Value = 6867.1665
$ws.Range("M75").Value = -5931.1665
$ws.Range("H78").Value = 13743.286
$ws.Range("I78").Value = 6867.1665
$ws.Range("K78").Value = 20601.4995
$ws.Range("M78").Value = -15921.4995
$ws.Range("H105").Value = 10072.818
$ws.Range("I105").Value = 13312.625
$ws.Range("J105").Value = 1433.3334
$ws.Range("K105").Value = 13312.625
$ws.Range("L105").Value = 1433.3334
$ws.Range("M105").Value = -11565.625
$ws.Range("N105").Value = -4927.3334
$ws.Range("H106").Value = 12937.8
$ws.Range("J106").Value = 12937.8
$ws.Range("L106").Value = 12937.8
$ws.Range("N106").Value = -15461.8
$ws.Range("H108").Value = 93999.5
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

# --- Sheet CRP: market price refresh ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1394.2354
$ws.Range("I16").Value = 1516.5385
$ws.Range("J16").Value = 996.75
$ws.Range("K16").Value = 1516.5385
$ws.Range("L16").Value = 996.75
$ws.Range("M16").Value = -1229.5385
$ws.Range("N16").Value = -1570.75
$ws.Range("H31").Value = 25004522
$ws.Range("I31").Value = 3504.724
$ws.Range("K31").Value = 3504.724
$ws.Range("M31").Value = -3209.724
$ws.Range("H34").Value = 25004522
$ws.Range("I34").Value = 3504.724
$ws.Range("K34").Value = 3504.724
$ws.Range("M34").Value = -3302.724
$ws.Range("H113").Value = 1394.2354
$ws.Range("I113").Value = 1516.5385
$ws.Range("J113").Value = 996.75
$ws.Range("K113").Value = 1516.5385
$ws.Range("L113").Value = 996.75
$ws.Range("M113").Value = 653.4614999999999
$ws.Range("N113").Value = -5336.75
$ws.Range("H132").Value = 120810.18
$ws.Range("I132").Value = 156195.77
$ws.Range("K132").Value = 468587.3099999999
$ws.Range("M132").Value = -466057.3099999999

# --- Sheet CUL: market price refresh ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 3105.1
$ws.Range("J39").Value = 4046.3076
$ws.Range("L39").Value = 12138.9228
$ws.Range("N39").Value = -12726.9228
$ws.Range("H55").Value = 2858244
$ws.Range("J55").Value = 1162.5807
$ws.Range("L55").Value = 3487.7421
$ws.Range("N55").Value = -3841.7421
$ws.Range("H69").Value = 474.5
$ws.Range("I69").Value = 474.5
$ws.Range("K69").Value = 1423.5
$ws.Range("M69").Value = -612.5
$ws.Range("H72").Value = 474.5
$ws.Range("I72").Value = 474.5
$ws.Range("K72").Value = 4270.5
$ws.Range("M72").Value = -214.5
$ws.Range("H80").Value = 4500
$ws.Range("J80").Value = 4500
$ws.Range("L80").Value = 13500
$ws.Range("N80").Value = -15372
$ws.Range("H83").Value = 4500
$ws.Range("J83").Value = 4500
$ws.Range("L83").Value = 40500
$ws.Range("N83").Value = -49860
$ws.Range("H121").Value = 2006
$ws.Range("I121").Value = 2010
$ws.Range("J121").Value = 2000
$ws.Range("K121").Value = 6030
$ws.Range("L121").Value = 6000
$ws.Range("M121").Value = -4720
$ws.Range("N121").Value = -8620
$ws.Range("H136").Value = 7250
$ws.Range("I136").Value = 4500
$ws.Range("J136").Value = 10000
$ws.Range("K136").Value = 13500
$ws.Range("L136").Value = 30000
$ws.Range("M136").Value = -8400
$ws.Range("N136").Value = -40200
$ws.Range("H138").Value = 1372.25
$ws.Range("I138").Value = 1346.7
$ws.Range("K138").Value = 4040.1
$ws.Range("M138").Value = 1099.9
$ws.Range("H141").Value = 7272.2085
$ws.Range("I141").Value = 5259.25
$ws.Range("J141").Value = 9285.166999999999
$ws.Range("K141").Value = 15777.75
$ws.Range("L141").Value = 27855.501
$ws.Range("M141").Value = -10597.75
$ws.Range("N141").Value = -38215.501

# --- Sheet GSM: market price refresh ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3930.0715
$ws.Range("I113").Value = 3177.3333
$ws.Range("K113").Value = 3177.3333
$ws.Range("M113").Value = -1007.3333
$ws.Range("H114").Value = 75000
$ws.Range("J114").Value = 75000
$ws.Range("L114").Value = 75000
$ws.Range("N114").Value = -83678
$ws.Range("H122").Value = 2058.9285
$ws.Range("I122").Value = 1818.75
$ws.Range("J122").Value = 3500
$ws.Range("K122").Value = 5456.25
$ws.Range("L122").Value = 10500
$ws.Range("M122").Value = -3006.25
$ws.Range("N122").Value = -15400

# --- Sheet LTW: market price refresh ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3648.75
$ws.Range("I22").Value = 695
$ws.Range("J22").Value = 4239.5
$ws.Range("K22").Value = 695
$ws.Range("L22").Value = 4239.5
$ws.Range("M22").Value = -400
$ws.Range("N22").Value = -4829.5
$ws.Range("H27").Value = 3648.75
$ws.Range("I27").Value = 695
$ws.Range("J27").Value = 4239.5
$ws.Range("K27").Value = 695
$ws.Range("L27").Value = 4239.5
$ws.Range("M27").Value = -588
$ws.Range("N27").Value = -4453.5
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H132").Value = 166668220
$ws.Range("I132").Value = 1875.9
$ws.Range("K132").Value = 5627.700000000001
$ws.Range("M132").Value = -3097.700000000001

# --- Sheet WVR: market price refresh ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 22752202
$ws.Range("I122").Value = 34519424
$ws.Range("K122").Value = 103558272
$ws.Range("M122").Value = -103555822
